$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.847390353679657
$ws.Range("B1").Value = 1.159092307090759
$ws.Range("C1").Value = 1.952577590942383
$ws.Range("D1").Value = 4.705146312713623
$ws.Range("E1").Value = 2.232603549957275
